$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header / value columns (D, E) ---------------------------------
# D1/D2 and E1/E2 are new columns added to the feedback sheet.
# Order matches shared-string insertion order from the source edit.
$ws.Range("E1").Value = "YE_Feedback_Knowledge"
$ws.Range("D2").Value = "Agreed_price"
$ws.Range("D1").Value = "YE_Estimate"
$ws.Range("E2").Value = "Yes"

# --- Column widths -------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 31.451822916666668
$ws.Columns.Item(3).ColumnWidth = 20.451822916666668
$ws.Columns.Item(4).ColumnWidth = 14.736979166666666
$ws.Columns.Item(5).ColumnWidth = 22.166666666666668

# --- Selection: Excel left the active cell on E1 when the file was saved.
$ws.Range("E1").Select()
